$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column A (shifts existing A:D -> B:E), preserving existing
# widths/values/styles on the shifted columns.
$ws.Columns.Item(1).Insert()

# New column A header + tab-name values.
$ws.Range("A1").Value = "TabName"
$ws.Range("A2").Value = "CasesTab"

# Updated Cypher queries (replace the old ones that used to live in A2/B2).
$casesQuery = @"
MATCH (ct:clinical_trial)<--(a:arm)<--(c:case)
    WHERE c.race = "ASIAN"
WITH DISTINCT c, a, ct
RETURN 
    COALESCE(c.case_id, '') AS ``Case ID``,
    COALESCE(ct.clinical_trial_designation, '') AS ``Trial Code``,
    COALESCE(a.arm_id, '') AS ``Arm``,
    COALESCE(a.arm_drug, '') AS ``Arm Treatment``,
    COALESCE(c.disease, '') AS ``Diagnosis``,
    COALESCE(c.gender, '') AS ``Gender``,
    COALESCE(c.race, '') AS ``Race``,
    COALESCE(c.ethnicity, '') AS ``Ethnicity``
"@

$statQuery = @"
MATCH (s:specimen)-->(c:case)-->(:arm)-->(ct:clinical_trial)
    WHERE c.race = "ASIAN"
OPTIONAL MATCH (f:file)-->(:sequencing_assay)-->(:nucleic_acid)-->(s)
RETURN 
    COUNT(DISTINCT f) AS number_of_files,
    COUNT(DISTINCT c.case_id) AS number_of_cases,
    COUNT(DISTINCT ct.clinical_trial_designation) AS number_of_trials
"@

$ws.Range("B2").Value = $casesQuery
$ws.Range("C2").Value = $statQuery

# C2 needs the same wrap-text style that B2 already carries.
$ws.Range("C2").WrapText = $true

# Row grew taller to fit the longer wrapped query text.
$ws.Rows.Item(2).RowHeight = 174

# Narrow auto-width column for the short TabName/CasesTab labels.
$ws.Columns.Item(1).ColumnWidth = 8

# Selection moved to B2 after the edit.
[void]$ws.Range("B2").Select()
